# "Generate Report for Handoff" - refresh the localization-status report.
#
# All rows that were "Ready for handoff" (or the failed-transform row) get
# their date/time column bumped to the new handoff run's timestamp:
#   - Overview sheet (column D, "Latest Handoff Date")
#   - zh-cn sheet     (column E, "Latest Handoff Datetime")
#   - de-de sheet     (column E, "Latest Handoff Datetime")

$wb = $excel.ActiveWorkbook

$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws = $wb.Worksheets.Item("Overview")
foreach ($r in $overviewRows) {
    $ws.Cells.Item($r, 4).Value = "2016-23-19 10:23:41"
}

$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $overviewRows) {
    $ws.Cells.Item($r, 5).Value = "2016-03-19 10:23:37"
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $overviewRows) {
    $ws.Cells.Item($r, 5).Value = "2016-03-19 10:23:41"
}
